# "setting for Powersafe changed and some fixes for Mega-sketch"
#
# - Mega 2560 sheet becomes the active/selected sheet (was Uno V3 before).
# - Its page setup gets a "fit to page" scale of 31%.
# - Its selection moves to B1:F71 (whole used range), with B1 active.
# - Uno V3 sheet is no longer the active/selected tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Mega 2560")

# Make "Mega 2560" the active sheet/tab (this also clears tabSelected on
# "Uno V3", which was previously the active tab).
$ws1.Activate() | Out-Null

# Powersafe / print setting change: scale the sheet to 31% and mark the
# sheet as using the "fit to page" page-setup mode.
$ws1.PageSetup.Zoom = 31
$ws1.PageSetup.FitToPagesWide = 1

# Select the sheet's whole used range, with B1 as the active cell.
$ws1.Range("B1:F71").Select() | Out-Null
